$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.772.80"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "'2.557.37"
$ws.Range("E3").Value = "  +5.86%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'574.34"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").Value = "'149.93"
$ws.Range("E6").Value = "  +8.14%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "'2.554.14"
$ws.Range("E9").Value = "  +5.81%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "'0.359"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "'28.26"
$ws.Range("E14").Value = "  +9.44%  "
$ws.Range("D15").Value = "'3.010.80"
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "'63.672.09"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "'2.549.47"
$ws.Range("E18").Value = "  +5.87%  "
$ws.Range("D19").Value = "'11.58"
$ws.Range("E19").Value = "  +4.56%  "
$ws.Range("D20").Value = "'343.10"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").Value = "'6.91"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'66.15"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "'1.59"
$ws.Range("E26").Value = "  +5.81%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'8.45"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("D30").Value = "'7.13"
$ws.Range("E30").Value = "  +12.74%  "
$ws.Range("D31").Value = "'0.0₃0837"
$ws.Range("E31").Value = "  +7.43%  "
$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = "  +4.71%  "
$ws.Range("D33").Value = "'177.30"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").Value = "'1.62"
$ws.Range("E34").Value = "  +14.16%  "
$ws.Range("D35").Value = "'423.35"
$ws.Range("E35").Value = "  +12.49%  "
$ws.Range("D36").Value = "'0.408"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("D37").Value = "'19.14"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +5.44%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").Value = "'156.59"
$ws.Range("E43").Value = "  +7.62%  "
$ws.Range("D44").Value = "'3.81"
$ws.Range("E44").Value = "  +4.20%  "
$ws.Range("D45").Value = "'21.05"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").Value = "'0.612"
$ws.Range("E46").Value = "  +4.42%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").Value = "'0.0969"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").Value = "'0.0233"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("D50").Value = "'18.83"
$ws.Range("E50").Value = "  +4.72%  "
$ws.Range("E51").Value = "  +11.26%  "
